$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gains two trailing
#    spaces, then four red-colored runs are appended:
#    "(This is a change – Ve" / "rsion for " / "main branch" / ")"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertAfter("  ")

$dash = [char]0x2013
$red = 255  # wdColorRed

$pos = $p1.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("(This is a change " + $dash + " Ve")
$ins.Font.Color = $red

$pos = $p1.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("rsion for ")
$ins.Font.Color = $red

$pos = $p1.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("main branch")
$ins.Font.Color = $red

$pos = $p1.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter(")")
$ins.Font.Color = $red

# ---------------------------------------------------------------------------
# 2) Add an extra blank paragraph right after "It will be treated as a binary
#    file by Git." (so there are two empty paragraphs there instead of one).
# ---------------------------------------------------------------------------
$gitPara = $d.Paragraphs(2)
$gitPara.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3) Move the phrase "even though it may, for the time being, anger them to
#    the point of cutting off the friendly hand that is stretched out to "
#    from the end of the run before the page break to the start of the run
#    after the page break (net text is unchanged; only the run split moves).
# ---------------------------------------------------------------------------
$moved = "even though it may, for the time being, anger them to the point of cutting off the friendly hand that is stretched out to "

$d.Content.Find.Execute(
    $moved + "help them.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "help them.", 2) | Out-Null

$d.Content.Find.Execute(
    "danger " + $moved,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "danger ", 2) | Out-Null

$d.Content.Find.Execute(
    "help them. People may laugh",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $moved + "help them. People may laugh", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) After the final (empty) "Normal (Web)" paragraph, add two more empty
#    paragraphs: a plain one, then one styled "larger" with the same
#    shading/spacing as the preceding "Normal (Web)" paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null
$last = $d.Paragraphs.Add()
$last.Range.Style = "larger"

# ---------------------------------------------------------------------------
# 5) Style "Normal (Web)" becomes semi-hidden.
# ---------------------------------------------------------------------------
$d.Styles("Normal (Web)").SemiHidden = $true
